$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1110.9756
$ws.Range("J17").Value = 1110.9756
$ws.Range("L17").Value = 3332.9268
$ws.Range("N17").Value = -3668.9268

$ws.Range("H112").Value = 6812.156
$ws.Range("J112").Value = 8112.3516
$ws.Range("L112").Value = 24337.0548
$ws.Range("N112").Value = -26553.0548

$ws.Range("H116").Value = 2193.8333
$ws.Range("I116").Value = 1758.3
$ws.Range("J116").Value = 2738.25
$ws.Range("K116").Value = 1758.3
$ws.Range("L116").Value = 2738.25
$ws.Range("M116").Value = 1683.7
$ws.Range("N116").Value = -9622.25

$ws.Range("H135").Value = 34483270
$ws.Range("I135").Value = 273.28
$ws.Range("J135").Value = 250002030
$ws.Range("K135").Value = 2459.52
$ws.Range("L135").Value = 2250018270
$ws.Range("M135").Value = 75.48000000000047
$ws.Range("N135").Value = -2250023340

$ws.Range("H137").Value = 1099.1714
$ws.Range("I137").Value = 849.1739
$ws.Range("J137").Value = 1578.3334
$ws.Range("K137").Value = 2547.5217
$ws.Range("L137").Value = 4735.0002
$ws.Range("M137").Value = 2.478299999999763
$ws.Range("N137").Value = -9835.0002

$ws.Range("H138").Value = 1102.6471
$ws.Range("I138").Value = 680.6889
$ws.Range("K138").Value = 2042.0667
$ws.Range("M138").Value = 3097.9333


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1566.6666
$ws.Range("I45").Value = 1614.2858
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 1614.2858
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -1237.2858
$ws.Range("N45").Value = -2154

$ws.Range("H61").Value = 29413072
$ws.Range("I61").Value = 38462720
$ws.Range("J61").Value = 1713.5
$ws.Range("K61").Value = 38462720
$ws.Range("L61").Value = 1713.5
$ws.Range("M61").Value = -38462508
$ws.Range("N61").Value = -2137.5

$ws.Range("H102").Value = 9805240
$ws.Range("I102").Value = 12821602
$ws.Range("J102").Value = 2065.5
$ws.Range("K102").Value = 12821602
$ws.Range("L102").Value = 2065.5
$ws.Range("M102").Value = -12819980
$ws.Range("N102").Value = -5309.5

$ws.Range("H136").Value = 29413072
$ws.Range("I136").Value = 38462720
$ws.Range("J136").Value = 1713.5
$ws.Range("K136").Value = 115388160
$ws.Range("L136").Value = 5140.5
$ws.Range("M136").Value = -115385610
$ws.Range("N136").Value = -10240.5


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 40001930
$ws.Range("I105").Value = 50001684
$ws.Range("K105").Value = 50001684
$ws.Range("M105").Value = -49999937

$ws.Range("H107").Value = 1143.6364
$ws.Range("I107").Value = 903.5
$ws.Range("K107").Value = 903.5
$ws.Range("M107").Value = 1016.5

$ws.Range("H134").Value = 4246.029
$ws.Range("J134").Value = 28750
$ws.Range("L134").Value = 86250
$ws.Range("N134").Value = -91320


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 62501200
$ws.Range("I16").Value = 66667844
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 66667844
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -66667557
$ws.Range("N16").Value = -2074

$ws.Range("H31").Value = 2388.1365
$ws.Range("I31").Value = 2484.647
$ws.Range("J31").Value = 2060
$ws.Range("K31").Value = 2484.647
$ws.Range("L31").Value = 2060
$ws.Range("M31").Value = -2189.647
$ws.Range("N31").Value = -2650

$ws.Range("H34").Value = 2388.1365
$ws.Range("I34").Value = 2484.647
$ws.Range("J34").Value = 2060
$ws.Range("K34").Value = 2484.647
$ws.Range("L34").Value = 2060
$ws.Range("M34").Value = -2282.647
$ws.Range("N34").Value = -2464

$ws.Range("H113").Value = 62501200
$ws.Range("I113").Value = 66667844
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 66667844
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -66665674
$ws.Range("N113").Value = -5840

$ws.Range("H129").Value = 43622
$ws.Range("J129").Value = 43622
$ws.Range("L129").Value = 43622
$ws.Range("N129").Value = -53622

$ws.Range("H135").Value = 30934.666
$ws.Range("J135").Value = 33551.5
$ws.Range("L135").Value = 33551.5
$ws.Range("N135").Value = -43691.5


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3394.8
$ws.Range("I94").Value = 3024
$ws.Range("J94").Value = 3487.5
$ws.Range("K94").Value = 9072
$ws.Range("L94").Value = 10462.5
$ws.Range("M94").Value = -8396
$ws.Range("N94").Value = -11814.5

$ws.Range("H109").Value = 78887.234
$ws.Range("I109").Value = 112170.445
$ws.Range("J109").Value = 4000
$ws.Range("K109").Value = 336511.335
$ws.Range("L109").Value = 12000
$ws.Range("M109").Value = -335471.335
$ws.Range("N109").Value = -14080

$ws.Range("H114").Value = 501.7
$ws.Range("I114").Value = 364.22223
$ws.Range("J114").Value = 614.1818
$ws.Range("K114").Value = 1092.66669
$ws.Range("L114").Value = 1842.5454
$ws.Range("M114").Value = 2161.33331
$ws.Range("N114").Value = -8350.545399999999

$ws.Range("H119").Value = 7397.6
$ws.Range("I119").Value = 3000
$ws.Range("J119").Value = 8074.154
$ws.Range("K119").Value = 9000
$ws.Range("L119").Value = 24222.462
$ws.Range("M119").Value = -4162
$ws.Range("N119").Value = -33898.462

$ws.Range("H129").Value = 13442518
$ws.Range("J129").Value = 3335265.5
$ws.Range("L129").Value = 10005796.5
$ws.Range("N129").Value = -10015796.5

$ws.Range("H131").Value = 15626165
$ws.Range("I131").Value = 90909570
$ws.Range("J131").Value = 1307.6038
$ws.Range("K131").Value = 272728710
$ws.Range("L131").Value = 3922.811400000001
$ws.Range("M131").Value = -272723670
$ws.Range("N131").Value = -14002.8114

$ws.Range("H132").Value = 1657.7
$ws.Range("I132").Value = 1446.5
$ws.Range("K132").Value = 13018.5
$ws.Range("M132").Value = -10488.5


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 13160472
$ws.Range("I102").Value = 19232318
$ws.Range("J102").Value = 4806.3335
$ws.Range("K102").Value = 19232318
$ws.Range("L102").Value = 4806.3335
$ws.Range("M102").Value = -19230696
$ws.Range("N102").Value = -8050.3335


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3577.9285
$ws.Range("I40").Value = 2198.7273
$ws.Range("K40").Value = 2198.7273
$ws.Range("M40").Value = -2062.7273


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 916.6667
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = $null

$ws.Range("H113").Value = 281.1154
$ws.Range("J113").Value = 343.5
$ws.Range("L113").Value = 1030.5
$ws.Range("N113").Value = -5370.5

$ws.Range("H136").Value = 371.5814
$ws.Range("I136").Value = 308.37143
$ws.Range("K136").Value = 925.11429
$ws.Range("M136").Value = 1624.88571

